$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 headers: Area / Atotal (new shared strings), duplicated Atotal/Qtotal in J11:K11
$ws.Range("G11").Value = "Area"
$ws.Range("H11").Value = "Atotal"
$ws.Range("J11").Value = "Atotal"
$ws.Range("K11").Value = "Qtotal"

# --- Row 12: Area + Atotal formulas, plus J12/K12 mirror formulas
$ws.Range("G12").Formula = "=(D12-0)*B12/100"
$ws.Range("H12").Formula = "=SUM(G12:G21)"
$ws.Range("J12").Formula = "=H12"
$ws.Range("K12").Formula = "=F12"

# --- Row 13: single Area formula
$ws.Range("G13").Formula = "=(D13-D12)*B13/100"

# --- Rows 14:25 share the Area formula (creates a shared formula G14:G25)
$ws.Range("G14:G25").Formula = "=(D14-D13)*B14/100"

# --- Selection matches the new focus cell used while editing
$ws.Range("J12:K12").Select

Write-Host "done"
